# Automated BRVM data refresh (Recommandations + Top_YTD) via GitHub Actions
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# ---- Recommandations sheet: refresh rows 2-47 (A:G) ----
$rec = New-Object 'object[,]' 46,7
$rec[0,0] = "SUCRIVOIRE"
$rec[0,1] = 0
$rec[0,2] = 4
$rec[0,3] = 3935
$rec[0,4] = 995
$rec[0,5] = "🟡 Observer"
$rec[0,6] = "➖ Neutre"
$rec[1,0] = "BRVM - SERVICES PUBLICS"
$rec[1,1] = 0
$rec[1,2] = 8
$rec[1,3] = 3329.56
$rec[1,4] = 112.02
$rec[1,5] = "🟡 Observer"
$rec[1,6] = "➖ Neutre"
$rec[2,0] = "SAFCA CI"
$rec[2,1] = 0
$rec[2,2] = 4
$rec[2,3] = 2795
$rec[2,4] = 695
$rec[2,5] = "🟡 Observer"
$rec[2,6] = "➖ Neutre"
$rec[3,0] = "CFAO MOTORS CI"
$rec[3,1] = 0
$rec[3,2] = 4
$rec[3,3] = 2700
$rec[3,4] = 675
$rec[3,5] = "🟡 Observer"
$rec[3,6] = "➖ Neutre"
$rec[4,0] = "BRVM - AUTRES SECTEURS"
$rec[4,1] = 0
$rec[4,2] = 4
$rec[4,3] = 2630.34
$rec[4,4] = 653.39
$rec[4,5] = "🟡 Observer"
$rec[4,6] = "➖ Neutre"
$rec[5,0] = "NEI-CEDA CI"
$rec[5,1] = 0
$rec[5,2] = 4
$rec[5,3] = 2385
$rec[5,4] = 595
$rec[5,5] = "🟡 Observer"
$rec[5,6] = "➖ Neutre"
$rec[6,0] = "UNIWAX CI"
$rec[6,1] = 0
$rec[6,2] = 4
$rec[6,3] = 2330
$rec[6,4] = 580
$rec[6,5] = "🟡 Observer"
$rec[6,6] = "➖ Neutre"
$rec[7,0] = "SETAO CI"
$rec[7,1] = 0
$rec[7,2] = 4
$rec[7,3] = 2210
$rec[7,4] = 555
$rec[7,5] = "🟡 Observer"
$rec[7,6] = "➖ Neutre"
$rec[8,0] = "AIR LIQUIDE CI"
$rec[8,1] = 0
$rec[8,2] = 4
$rec[8,3] = 2120
$rec[8,4] = 525
$rec[8,5] = "🟡 Observer"
$rec[8,6] = "➖ Neutre"
$rec[9,0] = "BRVM - DISTRIBUTION"
$rec[9,1] = 0
$rec[9,2] = 4
$rec[9,3] = 1468.45
$rec[9,4] = 367.04
$rec[9,5] = "🟡 Observer"
$rec[9,6] = "➖ Neutre"
$rec[10,0] = "BRVM - TRANSPORT"
$rec[10,1] = 0
$rec[10,2] = 4
$rec[10,3] = 1392.76
$rec[10,4] = 348.8
$rec[10,5] = "🟡 Observer"
$rec[10,6] = "➖ Neutre"
$rec[11,0] = "BRVM - AGRICULTURE"
$rec[11,1] = 0
$rec[11,2] = 4
$rec[11,3] = 1266.62
$rec[11,4] = 308.95
$rec[11,5] = "🟡 Observer"
$rec[11,6] = "➖ Neutre"
$rec[12,0] = "BRVM - INDUSTRIE"
$rec[12,1] = 0
$rec[12,2] = 4
$rec[12,3] = 800.05
$rec[12,4] = 202.29
$rec[12,5] = "🟡 Observer"
$rec[12,6] = "➖ Neutre"
$rec[13,0] = "BRVM-PRINCIPAL"
$rec[13,1] = 0
$rec[13,2] = 4
$rec[13,3] = 701.1900000000001
$rec[13,4] = 176.59
$rec[13,5] = "🟡 Observer"
$rec[13,6] = "➖ Neutre"
$rec[14,0] = "BRVM - CONSOMMATION DE BASE"
$rec[14,1] = 0
$rec[14,2] = 4
$rec[14,3] = 693.85
$rec[14,4] = 173.98
$rec[14,5] = "🟡 Observer"
$rec[14,6] = "➖ Neutre"
$rec[15,0] = "BRVM - INDUSTRIELS"
$rec[15,1] = 0
$rec[15,2] = 4
$rec[15,3] = 529.95
$rec[15,4] = 131.65
$rec[15,5] = "🟡 Observer"
$rec[15,6] = "➖ Neutre"
$rec[16,0] = "BRVM-PRESTIGE"
$rec[16,1] = 0
$rec[16,2] = 4
$rec[16,3] = 519.13
$rec[16,4] = 130.36
$rec[16,5] = "🟡 Observer"
$rec[16,6] = "➖ Neutre"
$rec[17,0] = "BRVM - FINANCES"
$rec[17,1] = 0
$rec[17,2] = 4
$rec[17,3] = 485.58
$rec[17,4] = 122.78
$rec[17,5] = "🟡 Observer"
$rec[17,6] = "➖ Neutre"
$rec[18,0] = "BRVM - SERVICES FINANCIERS"
$rec[18,1] = 0
$rec[18,2] = 4
$rec[18,3] = 477.23
$rec[18,4] = 120.67
$rec[18,5] = "🟡 Observer"
$rec[18,6] = "➖ Neutre"
$rec[19,0] = "BRVM - ENERGIE"
$rec[19,1] = 0
$rec[19,2] = 4
$rec[19,3] = 435.92
$rec[19,4] = 108.33
$rec[19,5] = "🟡 Observer"
$rec[19,6] = "➖ Neutre"
$rec[20,0] = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$rec[20,1] = 0
$rec[20,2] = 4
$rec[20,3] = 421.37
$rec[20,4] = 105.5
$rec[20,5] = "🟡 Observer"
$rec[20,6] = "➖ Neutre"
$rec[21,0] = "BRVM - TELECOMMUNICATIONS"
$rec[21,1] = 0
$rec[21,2] = 4
$rec[21,3] = 375.21
$rec[21,4] = 93.54000000000001
$rec[21,5] = "🟡 Observer"
$rec[21,6] = "➖ Neutre"
$rec[22,0] = "TRACTAFRIC MOTORS CI (PRSC)"
$rec[22,1] = 2
$rec[22,2] = 0
$rec[22,3] = 13.5
$rec[22,4] = 7.5
$rec[22,5] = "🟡 Observer"
$rec[22,6] = "➖ Neutre"
$rec[23,0] = "BANK OF AFRICA ML (BOAM)"
$rec[23,1] = 2
$rec[23,2] = 0
$rec[23,3] = 11.78
$rec[23,4] = 4.99
$rec[23,5] = "🟡 Observer"
$rec[23,6] = "➖ Neutre"
$rec[24,0] = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$rec[24,1] = 1
$rec[24,2] = 0
$rec[24,3] = 5.59
$rec[24,4] = 5.59
$rec[24,5] = "🟡 Observer"
$rec[24,6] = "➖ Neutre"
$rec[25,0] = "SICABLE CI (CABC)"
$rec[25,1] = 1
$rec[25,2] = 0
$rec[25,3] = 5.31
$rec[25,4] = 5.31
$rec[25,5] = "🟡 Observer"
$rec[25,6] = "➖ Neutre"
$rec[26,0] = "BANK OF AFRICA NG (BOAN)"
$rec[26,1] = 1
$rec[26,2] = 0
$rec[26,3] = 3.4
$rec[26,4] = 3.4
$rec[26,5] = "🟡 Observer"
$rec[26,6] = "➖ Neutre"
$rec[27,0] = "BANK OF AFRICA BN (BOAB)"
$rec[27,1] = 1
$rec[27,2] = 0
$rec[27,3] = 3.09
$rec[27,4] = 3.09
$rec[27,5] = "🟡 Observer"
$rec[27,6] = "➖ Neutre"
$rec[28,0] = "SONATEL SN (SNTS)"
$rec[28,1] = 1
$rec[28,2] = 1
$rec[28,3] = 2.08
$rec[28,4] = -1.92
$rec[28,5] = "🟡 Observer"
$rec[28,6] = "👀 À surveiller"
$rec[29,0] = "TOTALENERGIES MARKETING SN (TTLS)"
$rec[29,1] = 2
$rec[29,2] = 1
$rec[29,3] = 1.41
$rec[29,4] = -2.5
$rec[29,5] = "🟡 Observer"
$rec[29,6] = "👀 À surveiller"
$rec[30,0] = "VIVO ENERGY CI (SHEC)"
$rec[30,1] = 1
$rec[30,2] = 1
$rec[30,3] = 1.05
$rec[30,4] = -2.35
$rec[30,5] = "🟡 Observer"
$rec[30,6] = "👀 À surveiller"
$rec[31,0] = "ECOBANK COTE D''IVOIRE (ECOC)"
$rec[31,1] = 1
$rec[31,2] = 1
$rec[31,3] = 0.28
$rec[31,4] = 5.36
$rec[31,5] = "🟡 Observer"
$rec[31,6] = "👀 À surveiller"
$rec[32,0] = "TOTAL"
$rec[32,1] = 0
$rec[32,2] = 4
$rec[32,3] = 0
$rec[32,4] = 0
$rec[32,5] = "🟡 Observer"
$rec[32,6] = "➖ Neutre"
$rec[33,0] = "SODE CI (SDCC)"
$rec[33,1] = 1
$rec[33,2] = 1
$rec[33,3] = -0.6899999999999999
$rec[33,4] = 3.36
$rec[33,5] = "🟡 Observer"
$rec[33,6] = "👀 À surveiller"
$rec[34,0] = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$rec[34,1] = 1
$rec[34,2] = 1
$rec[34,3] = -0.71
$rec[34,4] = -2.11
$rec[34,5] = "🟡 Observer"
$rec[34,6] = "👀 À surveiller"
$rec[35,0] = "FILTISAC CI (FTSC)"
$rec[35,1] = 1
$rec[35,2] = 1
$rec[35,3] = -0.91
$rec[35,4] = -1.79
$rec[35,5] = "🟡 Observer"
$rec[35,6] = "👀 À surveiller"
$rec[36,0] = "BICI CI (BICC)"
$rec[36,1] = 0
$rec[36,2] = 1
$rec[36,3] = -2.54
$rec[36,4] = -2.54
$rec[36,5] = "🟡 Observer"
$rec[36,6] = "➖ Neutre"
$rec[37,0] = "SAFCA CI (SAFC)"
$rec[37,1] = 1
$rec[37,2] = 1
$rec[37,3] = -2.74
$rec[37,4] = -6.47
$rec[37,5] = "🟡 Observer"
$rec[37,6] = "👀 À surveiller"
$rec[38,0] = "BANK OF AFRICA SENEGAL (BOAS)"
$rec[38,1] = 0
$rec[38,2] = 1
$rec[38,3] = -3.37
$rec[38,4] = -3.37
$rec[38,5] = "🟡 Observer"
$rec[38,6] = "➖ Neutre"
$rec[39,0] = "ORANGE COTE D'IVOIRE (ORAC)"
$rec[39,1] = 0
$rec[39,2] = 1
$rec[39,3] = -3.38
$rec[39,4] = -3.38
$rec[39,5] = "🟡 Observer"
$rec[39,6] = "➖ Neutre"
$rec[40,0] = "UNIWAX CI (UNXC)"
$rec[40,1] = 1
$rec[40,2] = 2
$rec[40,3] = -3.84
$rec[40,4] = -4.35
$rec[40,5] = "🟡 Observer"
$rec[40,6] = "👀 À surveiller"
$rec[41,0] = "SOLIBRA CI (SLBC)"
$rec[41,1] = 1
$rec[41,2] = 2
$rec[41,3] = -4.85
$rec[41,4] = 7.5
$rec[41,5] = "🟡 Observer"
$rec[41,6] = "👀 À surveiller"
$rec[42,0] = "BERNABE CI (BNBC)"
$rec[42,1] = 1
$rec[42,2] = 2
$rec[42,3] = -5.09
$rec[42,4] = 4.74
$rec[42,5] = "🟡 Observer"
$rec[42,6] = "👀 À surveiller"
$rec[43,0] = "ECOBANK TRANS. INCORP. TG (ETIT)"
$rec[43,1] = 0
$rec[43,2] = 1
$rec[43,3] = -5.88
$rec[43,4] = -5.88
$rec[43,5] = "🟡 Observer"
$rec[43,6] = "➖ Neutre"
$rec[44,0] = "SETAO CI (STAC)"
$rec[44,1] = 0
$rec[44,2] = 1
$rec[44,3] = -6.09
$rec[44,4] = -6.09
$rec[44,5] = "🟡 Observer"
$rec[44,6] = "➖ Neutre"
$rec[45,0] = "SAPH CI (SPHC)"
$rec[45,1] = 0
$rec[45,2] = 1
$rec[45,3] = -7.46
$rec[45,4] = -7.46
$rec[45,5] = "🟡 Observer"
$rec[45,6] = "➖ Neutre"
$ws1.Range("A2:G47").Value = $rec

# Remove the 4 trailing rows that no longer exist in the refreshed dataset
$ws1.Rows("48:51").Delete()

# ---- Top_YTD sheet: refresh rows 2-11 (A:B) ----
$ytd = New-Object 'object[,]' 10,2
$ytd[0,0] = "BRVM - SERVICES PUBLICS"
$ytd[0,1] = 9072783.310000001
$ytd[1,0] = "SUCRIVOIRE"
$ytd[1,1] = 1379195.67
$ytd[2,0] = "SAFCA CI"
$ytd[2,1] = 406940
$ytd[3,0] = "CFAO MOTORS CI"
$ytd[3,1] = 360650.39
$ytd[4,0] = "BRVM - AUTRES SECTEURS"
$ytd[4,1] = 329279.09
$ytd[5,0] = "NEI-CEDA CI"
$ytd[5,1] = 234891.66
$ytd[6,0] = "UNIWAX CI"
$ytd[6,1] = 216823.4
$ytd[7,0] = "SETAO CI"
$ytd[7,1] = 180994.4
$ytd[8,0] = "AIR LIQUIDE CI"
$ytd[8,1] = 157419.69
$ytd[9,0] = "BRVM - DISTRIBUTION"
$ytd[9,1] = 47506.7
$ws2.Range("A2:B11").Value = $ytd

